# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450 -
# refresh the generated StructureDefinition spreadsheet (communication-retry-limt)
# with the new IG build's metadata + element Short/Definition text.

$wb = $excel.ActiveWorkbook

# ---- "Metadata" sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Build date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Range("B9").Value = "Alvearie Team"

# The old sheet had a duplicated "Contact / No display for ContactDetail"
# row pair; it is replaced by a single "Jurisdiction / United States of
# America" row, so drop the extra row first (rows shift up) and then set
# the remaining row's content.
$meta.Rows.Item(11).Delete()
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---- "Elements" sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) now carries the extension's own
# Short/Definition text instead of the generic placeholders.
$elements.Range("K2").Value = "Communication Retry Limit"
$elements.Range("L2").Value = "Retry limit configuration, may vary by medium and vendor"
